$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers need to be forced to text
# so they remain stored as strings (matching the source data which uses
# plain text price cells), rather than being auto-converted to numeric cells.
$textCells = @("D5", "D6", "D8", "D11", "D14", "D19", "D20", "D21", "D22", "D24", "D26", "D27", "D28", "D30", "D33", "D35", "D37", "D38", "D41", "D42", "D43", "D45", "D46", "D47", "D48", "D49", "D50")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Updated price (D) and volume-change (E) figures
$ws.Range("D2").Value = "63.491.06"
$ws.Range("E2").Value = "  +2.05%  "
$ws.Range("D3").Value = "2.545.45"
$ws.Range("E3").Value = "  +4.62%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "572.14"
$ws.Range("E5").Value = "  +2.78%  "
$ws.Range("D6").Value = "150.80"
$ws.Range("E6").Value = "  +8.59%  "
$ws.Range("D8").Value = "0.588"
$ws.Range("E8").Value = "  +0.98%  "
$ws.Range("D9").Value = "2.542.91"
$ws.Range("E9").Value = "  +4.64%  "
$ws.Range("E10").Value = "  +2.40%  "
$ws.Range("D11").Value = "5.76"
$ws.Range("E11").Value = "  +0.53%  "
$ws.Range("E12").Value = "  +1.05%  "
$ws.Range("E13").Value = "  +3.47%  "
$ws.Range("D14").Value = "28.44"
$ws.Range("E14").Value = "  +8.76%  "
$ws.Range("D15").Value = "3.002.22"
$ws.Range("E15").Value = "  +4.66%  "
$ws.Range("D16").Value = "63.380.31"
$ws.Range("E16").Value = "  +2.10%  "
$ws.Range("E17").Value = "  +1.36%  "
$ws.Range("D18").Value = "2.551.22"
$ws.Range("E18").Value = "  +4.97%  "
$ws.Range("D19").Value = "11.69"
$ws.Range("E19").Value = "  +4.29%  "
$ws.Range("D20").Value = "340.52"
$ws.Range("E20").Value = "  -1.58%  "
$ws.Range("D21").Value = "4.38"
$ws.Range("E21").Value = "  +4.28%  "
$ws.Range("D22").Value = "6.87"
$ws.Range("E22").Value = "  +1.26%  "
$ws.Range("E23").Value = "  +0.07%  "
$ws.Range("D24").Value = "66.22"
$ws.Range("E25").Value = "  -0.80%  "
$ws.Range("D26").Value = "1.60"
$ws.Range("E26").Value = "  +4.81%  "
$ws.Range("D27").Value = "1.51"
$ws.Range("E27").Value = "  +12.99%  "
$ws.Range("D28").Value = "8.48"
$ws.Range("E28").Value = "  +3.40%  "
$ws.Range("E29").Value = "  -0.07%  "
$ws.Range("D30").Value = "7.12"
$ws.Range("E30").Value = "  +11.91%  "
$ws.Range("D31").Value = "0.0₃0835"
$ws.Range("E31").Value = "  +6.65%  "
$ws.Range("E32").Value = "  +3.87%  "
$ws.Range("D33").Value = "177.93"
$ws.Range("E33").Value = "  +3.62%  "
$ws.Range("E34").Value = "  +9.25%  "
$ws.Range("D35").Value = "420.33"
$ws.Range("E35").Value = "  +14.88%  "
$ws.Range("E36").Value = "  +2.96%  "
$ws.Range("D37").Value = "19.19"
$ws.Range("E37").Value = "  +3.33%  "
$ws.Range("D38").Value = "4.45"
$ws.Range("E38").Value = "  -0.34%  "
$ws.Range("E40").Value = "  +4.48%  "
$ws.Range("D41").Value = "1.00"
$ws.Range("E41").Value = "  +0.03%  "
$ws.Range("D42").Value = "40.00"
$ws.Range("E42").Value = "  +2.00%  "
$ws.Range("D43").Value = "154.68"
$ws.Range("E43").Value = "  +5.77%  "
$ws.Range("E44").Value = "  +3.99%  "
$ws.Range("D45").Value = "21.16"
$ws.Range("E45").Value = "  +2.98%  "
$ws.Range("D46").Value = "0.612"
$ws.Range("E46").Value = "  +3.90%  "
$ws.Range("D47").Value = "0.0532"
$ws.Range("E47").Value = "  +2.86%  "
$ws.Range("D48").Value = "0.0242"
$ws.Range("E48").Value = "  +9.44%  "
$ws.Range("D49").Value = "0.0969"
$ws.Range("E49").Value = "  +1.64%  "
$ws.Range("D50").Value = "18.62"
$ws.Range("E50").Value = "  +4.26%  "
$ws.Range("E51").Value = "  +7.37%  "

# Restore default style on the cells we temporarily formatted as text,
# keeping their content as text while dropping the helper number format.
foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}
